$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (Player, Position, Team) for rows 2-18.
$players = @(
  @("Coby White", "PG,SG", "Chicago Bulls"),
  @("Norman Powell", "SG,SF", "LA Clippers"),
  @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
  @("Dorian Finney-Smith", "SF,PF,C", "Brooklyn Nets"),
  @("Guerschon Yabusele", "PF,C", "Philadelphia 76ers"),
  @("Walker Kessler", "C", "Utah Jazz"),
  @("Yves Missi", "C", "New Orleans Pelicans"),
  @("LeBron James", "SF,PF", "Los Angeles Lakers"),
  @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
  @("Jeremy Sochan", "SF,PF", "San Antonio Spurs"),
  @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
  @("Devin Booker", "PG,SG", "Phoenix Suns"),
  @("Jalen Brunson", "PG", "New York Knicks"),
  @("Alperen Sengün", "C", "Houston Rockets"),
  @("Trae Young", "PG", "Atlanta Hawks"),
  @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
  @("Kawhi Leonard", "SG,SF,PF", "LA Clippers")
)

for ($i = 0; $i -lt $players.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $players[$i][0]
  $ws.Cells.Item($row, 2).Value = $players[$i][1]
  $ws.Cells.Item($row, 3).Value = $players[$i][2]
}

# Original sheet had 18 data rows (19 total incl. header); the new roster
# only has 17 data rows, so remove the now-unused last row.
$ws.Rows.Item(19).Delete()
